$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Cash" row (original row 17); this shifts "1-5 years GILTS" up to row 17
$ws.Rows("17").Delete()

# Update the B/C numeric values for the remaining rows (2-17) to the new figures
$ws.Cells.Item(2, 2).Value2 = [double]"0.02998349699801747"
$ws.Cells.Item(2, 3).Value2 = [double]"0.02997980578320772"
$ws.Cells.Item(3, 2).Value2 = [double]"1.761589625444573e-19"
$ws.Cells.Item(3, 3).Value2 = [double]"3.862800410783512e-19"
$ws.Cells.Item(4, 2).Value2 = [double]"1.761589625444573e-19"
$ws.Cells.Item(4, 3).Value2 = [double]"1.62748231801186e-19"
$ws.Cells.Item(5, 2).Value2 = [double]"2.439454888092385e-19"
$ws.Cells.Item(5, 3).Value2 = [double]"1.62748231801186e-19"
$ws.Cells.Item(6, 2).Value2 = [double]"2.752045353098075e-19"
$ws.Cells.Item(6, 3).Value2 = [double]"2.846691045382254e-19"
$ws.Cells.Item(7, 2).Value2 = [double]"0.002899558148691509"
$ws.Cells.Item(7, 3).Value2 = [double]"0.002899779144328617"
$ws.Cells.Item(8, 2).Value2 = [double]"0.0219351662951515"
$ws.Cells.Item(8, 3).Value2 = [double]"0.0219363217695006"
$ws.Cells.Item(9, 2).Value2 = [double]"0.03294249137981119"
$ws.Cells.Item(9, 3).Value2 = [double]"0.0329428682566893"
$ws.Cells.Item(10, 2).Value2 = [double]"8.806444213787078e-20"
$ws.Cells.Item(10, 3).Value2 = [double]"1.423345522691127e-19"
$ws.Cells.Item(11, 2).Value2 = [double]"0.002411575005719938"
$ws.Cells.Item(11, 3).Value2 = [double]"0.002411783866887989"
$ws.Cells.Item(12, 2).Value2 = [double]"0.6949378640760929"
$ws.Cells.Item(12, 3).Value2 = [double]"0.6949474739127106"
$ws.Cells.Item(13, 2).Value2 = [double]"0.0008727300330695211"
$ws.Cells.Item(13, 3).Value2 = [double]"0.0008727562585592101"
$ws.Cells.Item(14, 2).Value2 = [double]"0.05277490247116246"
$ws.Cells.Item(14, 3).Value2 = [double]"0.05277251263754121"
$ws.Cells.Item(15, 2).Value2 = [double]"3.247273216924826e-19"
$ws.Cells.Item(15, 3).Value2 = [double]"4.95220678144552e-20"
$ws.Cells.Item(16, 2).Value2 = [double]"0.01684847136174504"
$ws.Cells.Item(16, 3).Value2 = [double]"0.01684995209365592"
$ws.Cells.Item(17, 2).Value2 = [double]"0.1443937442305386"
$ws.Cells.Item(17, 3).Value2 = [double]"0.1443867462769187"
